# Financials update for CANG_YR_FIN.xlsx
# Commit message: "Doing Updates for Financials"
# Updates the Income Statement / Balance Sheet / Cash Flow Statement figures
# on the "CANG" worksheet (columns D and E) to the latest reported numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CANG")

# --- Income Statement ---
$ws.Range("D8").Value  = 156200   # Total Revenue
$ws.Range("E8").Value  = 64500
$ws.Range("D9").Value  = 57300    # Cost of Revenue
$ws.Range("E9").Value  = 25200
$ws.Range("D10").Value = 98900    # Gross Profit
$ws.Range("E10").Value = 39200
$ws.Range("D12").Value = 2900     # Research Development
$ws.Range("D14").Value = -5800    # Non Recurring
$ws.Range("D17").Value = 86400    # Total Operating Expenses
$ws.Range("E17").Value = 37100
$ws.Range("D18").Value = 69800    # Operating Income or Loss
$ws.Range("E18").Value = 27400
$ws.Range("D21").Value = 71700    # Earnings Before Interest And Taxes
$ws.Range("E21").Value = "NA"
$ws.Range("D23").Value = 69500    # Income Before Tax
$ws.Range("E23").Value = 27700
$ws.Range("D24").Value = 17700    # Income Tax Expense
$ws.Range("E24").Value = 7900
$ws.Range("D26").Value = 51800    # Income After Tax
$ws.Range("E26").Value = 19800
$ws.Range("D27").Value = 25500    # Net Income From Continuing Ops
$ws.Range("E27").Value = 9600
$ws.Range("D33").Value = 25500    # Net Income
$ws.Range("E33").Value = 9600
$ws.Range("D35").Value = 25500    # Net Income Applicable To Common Shares
$ws.Range("E35").Value = 9600

# --- Balance Sheet ---
$ws.Range("D41").Value = 119200   # Cash And Cash Equivalents
$ws.Range("E41").Value = 6700
$ws.Range("D42").Value = 9300     # Short Term Investments
$ws.Range("E42").Value = 15700
$ws.Range("D43").Value = 23000    # Net Receivables
$ws.Range("E43").Value = 20600
$ws.Range("D45").Value = 13000    # Other Current Assets
$ws.Range("D46").Value = 164500   # Total Current Assets
$ws.Range("E46").Value = 44800
$ws.Range("D47").Value = 71900    # Long Term Investments
$ws.Range("E47").Value = 51000
$ws.Range("D49").Value = 300      # Goodwill
$ws.Range("D52").Value = 58300    # Other Assets
$ws.Range("E52").Value = 9500
$ws.Range("D54").Value = 296400   # Total Assets
$ws.Range("E54").Value = 106100
$ws.Range("D59").Value = 76700    # Accounts Payable
$ws.Range("E59").Value = 41800
$ws.Range("D60").Value = 78100    # Short/Current Long Term Debt
$ws.Range("E60").Value = 42600
$ws.Range("D61").Value = 26000    # Other Current Liabilities
$ws.Range("E61").Value = 28100
$ws.Range("D62").Value = 5300     # Total Current Liabilities
$ws.Range("E62").Value = 4100
$ws.Range("D66").Value = 113200   # Total Liabilities
$ws.Range("E66").Value = 77000
$ws.Range("D70").Value = 585000   # Common Stock
$ws.Range("E70").Value = 585000
$ws.Range("D72").Value = -402400  # Retained Earnings
$ws.Range("E72").Value = -556600
$ws.Range("D76").Value = -401800  # Total Stockholder Equity
$ws.Range("E76").Value = -555900

# --- Cash Flow Statement ---
$ws.Range("D81").Value  = 25500   # Net Income
$ws.Range("E81").Value  = 9600
$ws.Range("D83").Value  = 200     # Depreciation
$ws.Range("E83").Value  = "NA"
$ws.Range("D89").Value  = 12300   # Total Cash Flow From Operating Activities
$ws.Range("E89").Value  = "NA"
$ws.Range("D91").Value  = -400    # Capital Expenditures
$ws.Range("E91").Value  = "NA"
$ws.Range("D94").Value  = -40400  # Total Cash Flows From Investing Activities
$ws.Range("E94").Value  = "NA"
$ws.Range("D100").Value = 27800   # Total Cash Flows From Financing Activities
$ws.Range("E100").Value = "NA"
$ws.Range("D101").Value = 0       # Effect Of Exchange Rate Changes
$ws.Range("E101").Value = "NA"
$ws.Range("D102").Value = -300    # Change In Cash and Cash Equivalents
$ws.Range("E102").Value = "NA"
